# update tables with feedback from SC
$d = $word.ActiveDocument

# --- 1. Add a superscript "*" after "Model" in the header row ---
$rng = $d.Content
$rng.Find.Execute("Model", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)   # wdCollapseEnd
$rng.InsertAfter("*")
$rng.Font.Superscript = $true

# --- 2. Simple text replacements ---
$d.Content.Find.Execute("Patients", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Inpatient admissions", 2) | Out-Null

$d.Content.Find.Execute("Patient days", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Patient days (truncated at 14 days)", 2) | Out-Null

$d.Content.Find.Execute("965,391", $true, $false, $false, $false, $false,
                         $true, 1, $false, "935,735", 2) | Out-Null

$d.Content.Find.Execute("1,527,155", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1,244,944", 2) | Out-Null

$d.Content.Find.Execute("1,819,666", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1,532,182", 2) | Out-Null

$d.Content.Find.Execute("1,798,044", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1,501,363", 2) | Out-Null

$d.Content.Find.Execute("1,999,848", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1,706,308", 2) | Out-Null

$d.Content.Find.Execute("Full", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Final", 2) | Out-Null

$d.Content.Find.Execute("2,027,526", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1,730,133", 2) | Out-Null

# --- 3. Append a new footnote row to the bottom of the table ---
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$cells = $newRow.Cells
$cells.Item(1).Merge($cells.Item(6)) | Out-Null
$cell = $newRow.Cells.Item(1)

$bTop = $cell.Borders.Item(-1)      # wdBorderTop
$bTop.Color = 6710886               # 666666
$bTop.LineWidth = 6                 # -> sz 12
$bTop.LineStyle = 1                 # wdLineStyleSingle

$bLeft = $cell.Borders.Item(-2)     # wdBorderLeft
$bLeft.Color = 16777215             # FFFFFF
$bLeft.LineWidth = 0
$bLeft.LineStyle = 0                # wdLineStyleNone

$bBottom = $cell.Borders.Item(-3)   # wdBorderBottom
$bBottom.Color = 16777215           # FFFFFF
$bBottom.LineWidth = 0
$bBottom.LineStyle = 0

$bRight = $cell.Borders.Item(-4)    # wdBorderRight
$bRight.Color = 16777215            # FFFFFF
$bRight.LineWidth = 0
$bRight.LineStyle = 0

$cellRange = $cell.Range
$footStart = $cellRange.Start
$footnoteText = "Model represents the cross-validation fold models and the final model fit with all patient data. The fold models are those fit during internal-external cross-validation and incorporate all patient data except for the associated hospital of the same number. For example, the 'Fold: 1' model was fit using patient data from hospitals 2 to 5, with hospital 1 being the validation set."
$cellRange.Text = "*" + $footnoteText

$starRange = $d.Range($footStart, $footStart + 1)
$starRange.Font.Superscript = $true
